$wb = $excel.ActiveWorkbook

# --- Sheet1: Closing_Price ---
$ws1 = $wb.Worksheets.Item("Closing_Price")

# Update header B1 from "bitcoin" to "Trace_2"
$ws1.Range("B1").Value = "Trace_2"

# Pre-format the new rows (A3636:A3869) by copying the style of A3635 (date format)
$ws1.Range("A3635").Copy($ws1.Range("A3636:A3869"))

# Fill in new date/value rows (3636-3869)
$ws1.Cells.Item(3636, 1).Value = 45025
$ws1.Cells.Item(3636, 2).Value = 28351.23699385011
$ws1.Cells.Item(3637, 1).Value = 45026
$ws1.Cells.Item(3637, 2).Value = 29657.97413687356
$ws1.Cells.Item(3638, 1).Value = 45027
$ws1.Cells.Item(3638, 2).Value = 30260.93610940865
$ws1.Cells.Item(3639, 1).Value = 45028
$ws1.Cells.Item(3639, 2).Value = 29904.13869471891
$ws1.Cells.Item(3640, 1).Value = 45029
$ws1.Cells.Item(3640, 2).Value = 30405.02731278115
$ws1.Cells.Item(3641, 1).Value = 45030
$ws1.Cells.Item(3641, 2).Value = 30468.40870059078
$ws1.Cells.Item(3642, 1).Value = 45031
$ws1.Cells.Item(3642, 2).Value = 30312.16187965924
$ws1.Cells.Item(3643, 1).Value = 45032
$ws1.Cells.Item(3643, 2).Value = 30304.80751478584
$ws1.Cells.Item(3644, 1).Value = 45033
$ws1.Cells.Item(3644, 2).Value = 29467.45982926034
$ws1.Cells.Item(3645, 1).Value = 45034
$ws1.Cells.Item(3645, 2).Value = 30365.90416754141
$ws1.Cells.Item(3646, 1).Value = 45035
$ws1.Cells.Item(3646, 2).Value = 28833.2175012211
$ws1.Cells.Item(3647, 1).Value = 45036
$ws1.Cells.Item(3647, 2).Value = 28255.57824866478
$ws1.Cells.Item(3648, 1).Value = 45037
$ws1.Cells.Item(3648, 2).Value = 27300.15712851557
$ws1.Cells.Item(3649, 1).Value = 45038
$ws1.Cells.Item(3649, 2).Value = 27861.64066347655
$ws1.Cells.Item(3650, 1).Value = 45039
$ws1.Cells.Item(3650, 2).Value = 27606.57834765314
$ws1.Cells.Item(3651, 1).Value = 45040
$ws1.Cells.Item(3651, 2).Value = 27511.63568216131
$ws1.Cells.Item(3652, 1).Value = 45041
$ws1.Cells.Item(3652, 2).Value = 28351.21824837301
$ws1.Cells.Item(3653, 1).Value = 45042
$ws1.Cells.Item(3653, 2).Value = 28352.19132059668
$ws1.Cells.Item(3654, 1).Value = 45043
$ws1.Cells.Item(3654, 2).Value = 29483.52170499186
$ws1.Cells.Item(3655, 1).Value = 45044
$ws1.Cells.Item(3655, 2).Value = 29339.99499754935
$ws1.Cells.Item(3656, 1).Value = 45045
$ws1.Cells.Item(3656, 2).Value = 29217.94404680991
$ws1.Cells.Item(3657, 1).Value = 45046
$ws1.Cells.Item(3657, 2).Value = 29362.05621362539
$ws1.Cells.Item(3658, 1).Value = 45047
$ws1.Cells.Item(3658, 2).Value = 28125.50115563682
$ws1.Cells.Item(3659, 1).Value = 45048
$ws1.Cells.Item(3659, 2).Value = 28654.39013332672
$ws1.Cells.Item(3660, 1).Value = 45049
$ws1.Cells.Item(3660, 2).Value = 28988.32099624927
$ws1.Cells.Item(3661, 1).Value = 45050
$ws1.Cells.Item(3661, 2).Value = 28846.46145860806
$ws1.Cells.Item(3662, 1).Value = 45051
$ws1.Cells.Item(3662, 2).Value = 29520.3222689705
$ws1.Cells.Item(3663, 1).Value = 45052
$ws1.Cells.Item(3663, 2).Value = 28887.74104552337
$ws1.Cells.Item(3664, 1).Value = 45053
$ws1.Cells.Item(3664, 2).Value = 28611.43919761457
$ws1.Cells.Item(3665, 1).Value = 45054
$ws1.Cells.Item(3665, 2).Value = 27696.76078562556
$ws1.Cells.Item(3666, 1).Value = 45055
$ws1.Cells.Item(3666, 2).Value = 27607.39168091902
$ws1.Cells.Item(3667, 1).Value = 45056
$ws1.Cells.Item(3667, 2).Value = 27639.73356593586
$ws1.Cells.Item(3668, 1).Value = 45057
$ws1.Cells.Item(3668, 2).Value = 27024.76572929978
$ws1.Cells.Item(3669, 1).Value = 45058
$ws1.Cells.Item(3669, 2).Value = 26787.69039591833
$ws1.Cells.Item(3670, 1).Value = 45059
$ws1.Cells.Item(3670, 2).Value = 26798.1262714747
$ws1.Cells.Item(3671, 1).Value = 45060
$ws1.Cells.Item(3671, 2).Value = 26911.80527321389
$ws1.Cells.Item(3672, 1).Value = 45061
$ws1.Cells.Item(3672, 2).Value = 27227.79342255879
$ws1.Cells.Item(3673, 1).Value = 45062
$ws1.Cells.Item(3673, 2).Value = 27022.71317392115
$ws1.Cells.Item(3674, 1).Value = 45063
$ws1.Cells.Item(3674, 2).Value = 27389.9715114422
$ws1.Cells.Item(3675, 1).Value = 45064
$ws1.Cells.Item(3675, 2).Value = 26842.95249471792
$ws1.Cells.Item(3676, 1).Value = 45065
$ws1.Cells.Item(3676, 2).Value = 26884.37105958995
$ws1.Cells.Item(3677, 1).Value = 45066
$ws1.Cells.Item(3677, 2).Value = 27093.79110149735
$ws1.Cells.Item(3678, 1).Value = 45067
$ws1.Cells.Item(3678, 2).Value = 26773.83024366478
$ws1.Cells.Item(3679, 1).Value = 45068
$ws1.Cells.Item(3679, 2).Value = 26869.68602179078
$ws1.Cells.Item(3680, 1).Value = 45069
$ws1.Cells.Item(3680, 2).Value = 27222.93774660385
$ws1.Cells.Item(3681, 1).Value = 45070
$ws1.Cells.Item(3681, 2).Value = 26338.94894699869
$ws1.Cells.Item(3682, 1).Value = 45071
$ws1.Cells.Item(3682, 2).Value = 26475.60790253427
$ws1.Cells.Item(3683, 1).Value = 45072
$ws1.Cells.Item(3683, 2).Value = 26717.98755357178
$ws1.Cells.Item(3684, 1).Value = 45073
$ws1.Cells.Item(3684, 2).Value = 26848.23993959246
$ws1.Cells.Item(3685, 1).Value = 45074
$ws1.Cells.Item(3685, 2).Value = 28110.31346760928
$ws1.Cells.Item(3686, 1).Value = 45075
$ws1.Cells.Item(3686, 2).Value = 27759.74666392068
$ws1.Cells.Item(3687, 1).Value = 45076
$ws1.Cells.Item(3687, 2).Value = 27713.90854752951
$ws1.Cells.Item(3688, 1).Value = 45077
$ws1.Cells.Item(3688, 2).Value = 27245.47446517972
$ws1.Cells.Item(3689, 1).Value = 45078
$ws1.Cells.Item(3689, 2).Value = 26824.10149951678
$ws1.Cells.Item(3690, 1).Value = 45079
$ws1.Cells.Item(3690, 2).Value = 27247.74008145703
$ws1.Cells.Item(3691, 1).Value = 45080
$ws1.Cells.Item(3691, 2).Value = 27066.37132852623
$ws1.Cells.Item(3692, 1).Value = 45081
$ws1.Cells.Item(3692, 2).Value = 27315.44758748768
$ws1.Cells.Item(3693, 1).Value = 45082
$ws1.Cells.Item(3693, 2).Value = 25792.6160902334
$ws1.Cells.Item(3694, 1).Value = 45083
$ws1.Cells.Item(3694, 2).Value = 27216.61556469773
$ws1.Cells.Item(3695, 1).Value = 45084
$ws1.Cells.Item(3695, 2).Value = 26346.24454513763
$ws1.Cells.Item(3696, 1).Value = 45085
$ws1.Cells.Item(3696, 2).Value = 26507.90989220298
$ws1.Cells.Item(3697, 1).Value = 45086
$ws1.Cells.Item(3697, 2).Value = 26469.58168400703
$ws1.Cells.Item(3698, 1).Value = 45087
$ws1.Cells.Item(3698, 2).Value = 25858.12289218562
$ws1.Cells.Item(3699, 1).Value = 45088
$ws1.Cells.Item(3699, 2).Value = 25916.57989637981
$ws1.Cells.Item(3700, 1).Value = 45089
$ws1.Cells.Item(3700, 2).Value = 25910.36274291946
$ws1.Cells.Item(3701, 1).Value = 45090
$ws1.Cells.Item(3701, 2).Value = 25872.20645879509
$ws1.Cells.Item(3702, 1).Value = 45091
$ws1.Cells.Item(3702, 2).Value = 25107.75469588059
$ws1.Cells.Item(3703, 1).Value = 45092
$ws1.Cells.Item(3703, 2).Value = 25564.59963288714
$ws1.Cells.Item(3704, 1).Value = 45093
$ws1.Cells.Item(3704, 2).Value = 26327.3256694539
$ws1.Cells.Item(3705, 1).Value = 45094
$ws1.Cells.Item(3705, 2).Value = 26501.04444223367
$ws1.Cells.Item(3706, 1).Value = 45095
$ws1.Cells.Item(3706, 2).Value = 26333.09252966437
$ws1.Cells.Item(3707, 1).Value = 45096
$ws1.Cells.Item(3707, 2).Value = 26779.38746152284
$ws1.Cells.Item(3708, 1).Value = 45097
$ws1.Cells.Item(3708, 2).Value = 28330.89141980751
$ws1.Cells.Item(3709, 1).Value = 45098
$ws1.Cells.Item(3709, 2).Value = 30101.76481948869
$ws1.Cells.Item(3710, 1).Value = 45099
$ws1.Cells.Item(3710, 2).Value = 29935.63210674956
$ws1.Cells.Item(3711, 1).Value = 45100
$ws1.Cells.Item(3711, 2).Value = 30629.24435333575
$ws1.Cells.Item(3712, 1).Value = 45101
$ws1.Cells.Item(3712, 2).Value = 30537.81666463479
$ws1.Cells.Item(3713, 1).Value = 45102
$ws1.Cells.Item(3713, 2).Value = 30454.75697441562
$ws1.Cells.Item(3714, 1).Value = 45103
$ws1.Cells.Item(3714, 2).Value = 30285.50506407074
$ws1.Cells.Item(3715, 1).Value = 45104
$ws1.Cells.Item(3715, 2).Value = 30693.54635606798
$ws1.Cells.Item(3716, 1).Value = 45105
$ws1.Cells.Item(3716, 2).Value = 30083.47724526382
$ws1.Cells.Item(3717, 1).Value = 45106
$ws1.Cells.Item(3717, 2).Value = 30466.61045689353
$ws1.Cells.Item(3718, 1).Value = 45107
$ws1.Cells.Item(3718, 2).Value = 30480.78148194446
$ws1.Cells.Item(3719, 1).Value = 45108
$ws1.Cells.Item(3719, 2).Value = 30583.61959969754
$ws1.Cells.Item(3720, 1).Value = 45109
$ws1.Cells.Item(3720, 2).Value = 30571.88707287145
$ws1.Cells.Item(3721, 1).Value = 45110
$ws1.Cells.Item(3721, 2).Value = 31134.71313350841
$ws1.Cells.Item(3722, 1).Value = 45111
$ws1.Cells.Item(3722, 2).Value = 30775.6240223753
$ws1.Cells.Item(3723, 1).Value = 45112
$ws1.Cells.Item(3723, 2).Value = 30484.79377442352
$ws1.Cells.Item(3724, 1).Value = 45113
$ws1.Cells.Item(3724, 2).Value = 29990.05515122025
$ws1.Cells.Item(3725, 1).Value = 45114
$ws1.Cells.Item(3725, 2).Value = 30315.26540235163
$ws1.Cells.Item(3726, 1).Value = 45115
$ws1.Cells.Item(3726, 2).Value = 30264.19734525722
$ws1.Cells.Item(3727, 1).Value = 45116
$ws1.Cells.Item(3727, 2).Value = 30169.87065747925
$ws1.Cells.Item(3728, 1).Value = 45117
$ws1.Cells.Item(3728, 2).Value = 30394.28114637557
$ws1.Cells.Item(3729, 1).Value = 45118
$ws1.Cells.Item(3729, 2).Value = 30620.81424697801
$ws1.Cells.Item(3730, 1).Value = 45119
$ws1.Cells.Item(3730, 2).Value = 30407.41631435991
$ws1.Cells.Item(3731, 1).Value = 45120
$ws1.Cells.Item(3731, 2).Value = 31446.0149710716
$ws1.Cells.Item(3732, 1).Value = 45121
$ws1.Cells.Item(3732, 2).Value = 30311.82780813052
$ws1.Cells.Item(3733, 1).Value = 45122
$ws1.Cells.Item(3733, 2).Value = 30297.97140012553
$ws1.Cells.Item(3734, 1).Value = 45123
$ws1.Cells.Item(3734, 2).Value = 30237.3347851461
$ws1.Cells.Item(3735, 1).Value = 45124
$ws1.Cells.Item(3735, 2).Value = 30147.86007921152
$ws1.Cells.Item(3736, 1).Value = 45125
$ws1.Cells.Item(3736, 2).Value = 29848.03227250528
$ws1.Cells.Item(3737, 1).Value = 45126
$ws1.Cells.Item(3737, 2).Value = 29919.87483825516
$ws1.Cells.Item(3738, 1).Value = 45127
$ws1.Cells.Item(3738, 2).Value = 29787.04153035569
$ws1.Cells.Item(3739, 1).Value = 45128
$ws1.Cells.Item(3739, 2).Value = 29914.68272835899
$ws1.Cells.Item(3740, 1).Value = 45129
$ws1.Cells.Item(3740, 2).Value = 29710.15593599177
$ws1.Cells.Item(3741, 1).Value = 45130
$ws1.Cells.Item(3741, 2).Value = 30057.95671871459
$ws1.Cells.Item(3742, 1).Value = 45131
$ws1.Cells.Item(3742, 2).Value = 29184.90191305769
$ws1.Cells.Item(3743, 1).Value = 45132
$ws1.Cells.Item(3743, 2).Value = 29222.97463038353
$ws1.Cells.Item(3744, 1).Value = 45133
$ws1.Cells.Item(3744, 2).Value = 29363.67290797862
$ws1.Cells.Item(3745, 1).Value = 45134
$ws1.Cells.Item(3745, 2).Value = 29200.24400129131
$ws1.Cells.Item(3746, 1).Value = 45135
$ws1.Cells.Item(3746, 2).Value = 29314.09738709836
$ws1.Cells.Item(3747, 1).Value = 45136
$ws1.Cells.Item(3747, 2).Value = 29356.74774591667
$ws1.Cells.Item(3748, 1).Value = 45137
$ws1.Cells.Item(3748, 2).Value = 29277.75581027272
$ws1.Cells.Item(3749, 1).Value = 45138
$ws1.Cells.Item(3749, 2).Value = 29233.13636558044
$ws1.Cells.Item(3750, 1).Value = 45139
$ws1.Cells.Item(3750, 2).Value = 29537.10974692743
$ws1.Cells.Item(3751, 1).Value = 45140
$ws1.Cells.Item(3751, 2).Value = 29146.63627084363
$ws1.Cells.Item(3752, 1).Value = 45141
$ws1.Cells.Item(3752, 2).Value = 29176.15059658541
$ws1.Cells.Item(3753, 1).Value = 45142
$ws1.Cells.Item(3753, 2).Value = 29087.93808600763
$ws1.Cells.Item(3754, 1).Value = 45143
$ws1.Cells.Item(3754, 2).Value = 29046.8004570232
$ws1.Cells.Item(3755, 1).Value = 45144
$ws1.Cells.Item(3755, 2).Value = 29044.20437557756
$ws1.Cells.Item(3756, 1).Value = 45145
$ws1.Cells.Item(3756, 2).Value = 29178.14775187308
$ws1.Cells.Item(3757, 1).Value = 45146
$ws1.Cells.Item(3757, 2).Value = 29779.56167125089
$ws1.Cells.Item(3758, 1).Value = 45147
$ws1.Cells.Item(3758, 2).Value = 29585.48807737921
$ws1.Cells.Item(3759, 1).Value = 45148
$ws1.Cells.Item(3759, 2).Value = 29423.81891597763
$ws1.Cells.Item(3760, 1).Value = 45149
$ws1.Cells.Item(3760, 2).Value = 29396.8479714205
$ws1.Cells.Item(3761, 1).Value = 45150
$ws1.Cells.Item(3761, 2).Value = 29412.1422745841
$ws1.Cells.Item(3762, 1).Value = 45151
$ws1.Cells.Item(3762, 2).Value = 29284.96971374381
$ws1.Cells.Item(3763, 1).Value = 45152
$ws1.Cells.Item(3763, 2).Value = 29400.58680419105
$ws1.Cells.Item(3764, 1).Value = 45153
$ws1.Cells.Item(3764, 2).Value = 29170.49039706061
$ws1.Cells.Item(3765, 1).Value = 45154
$ws1.Cells.Item(3765, 2).Value = 28754.19702111257
$ws1.Cells.Item(3766, 1).Value = 45155
$ws1.Cells.Item(3766, 2).Value = 26501.58726991974
$ws1.Cells.Item(3767, 1).Value = 45156
$ws1.Cells.Item(3767, 2).Value = 26042.838256849
$ws1.Cells.Item(3768, 1).Value = 45157
$ws1.Cells.Item(3768, 2).Value = 26104.7786307979
$ws1.Cells.Item(3769, 1).Value = 45158
$ws1.Cells.Item(3769, 2).Value = 26160.70041981821
$ws1.Cells.Item(3770, 1).Value = 45159
$ws1.Cells.Item(3770, 2).Value = 26119.00595859522
$ws1.Cells.Item(3771, 1).Value = 45160
$ws1.Cells.Item(3771, 2).Value = 26033.67540658584
$ws1.Cells.Item(3772, 1).Value = 45161
$ws1.Cells.Item(3772, 2).Value = 26450.00862461353
$ws1.Cells.Item(3773, 1).Value = 45162
$ws1.Cells.Item(3773, 2).Value = 26134.934254516
$ws1.Cells.Item(3774, 1).Value = 45163
$ws1.Cells.Item(3774, 2).Value = 26044.03620928913
$ws1.Cells.Item(3775, 1).Value = 45164
$ws1.Cells.Item(3775, 2).Value = 26002.01571731937
$ws1.Cells.Item(3776, 1).Value = 45165
$ws1.Cells.Item(3776, 2).Value = 26082.71717821853
$ws1.Cells.Item(3777, 1).Value = 45166
$ws1.Cells.Item(3777, 2).Value = 26109.2698239287
$ws1.Cells.Item(3778, 1).Value = 45167
$ws1.Cells.Item(3778, 2).Value = 27730.50373550554
$ws1.Cells.Item(3779, 1).Value = 45168
$ws1.Cells.Item(3779, 2).Value = 27297.2614384804
$ws1.Cells.Item(3780, 1).Value = 45169
$ws1.Cells.Item(3780, 2).Value = 25927.41700577929
$ws1.Cells.Item(3781, 1).Value = 45170
$ws1.Cells.Item(3781, 2).Value = 25812.33049442542
$ws1.Cells.Item(3782, 1).Value = 45171
$ws1.Cells.Item(3782, 2).Value = 25853.65684277757
$ws1.Cells.Item(3783, 1).Value = 45172
$ws1.Cells.Item(3783, 2).Value = 25959.59631146345
$ws1.Cells.Item(3784, 1).Value = 45173
$ws1.Cells.Item(3784, 2).Value = 25829.36477294132
$ws1.Cells.Item(3785, 1).Value = 45174
$ws1.Cells.Item(3785, 2).Value = 25784.41351983038
$ws1.Cells.Item(3786, 1).Value = 45175
$ws1.Cells.Item(3786, 2).Value = 25752.95841858941
$ws1.Cells.Item(3787, 1).Value = 45176
$ws1.Cells.Item(3787, 2).Value = 26192.33343309057
$ws1.Cells.Item(3788, 1).Value = 45177
$ws1.Cells.Item(3788, 2).Value = 25907.22813724973
$ws1.Cells.Item(3789, 1).Value = 45178
$ws1.Cells.Item(3789, 2).Value = 25889.3250094629
$ws1.Cells.Item(3790, 1).Value = 45179
$ws1.Cells.Item(3790, 2).Value = 25834.58009856291
$ws1.Cells.Item(3791, 1).Value = 45180
$ws1.Cells.Item(3791, 2).Value = 25133.30310656653
$ws1.Cells.Item(3792, 1).Value = 45181
$ws1.Cells.Item(3792, 2).Value = 25866.80666315808
$ws1.Cells.Item(3793, 1).Value = 45182
$ws1.Cells.Item(3793, 2).Value = 26223.43607566095
$ws1.Cells.Item(3794, 1).Value = 45183
$ws1.Cells.Item(3794, 2).Value = 26531.39556626326
$ws1.Cells.Item(3795, 1).Value = 45184
$ws1.Cells.Item(3795, 2).Value = 26634.63102500059
$ws1.Cells.Item(3796, 1).Value = 45185
$ws1.Cells.Item(3796, 2).Value = 26557.76869199465
$ws1.Cells.Item(3797, 1).Value = 45186
$ws1.Cells.Item(3797, 2).Value = 26520.98825478389
$ws1.Cells.Item(3798, 1).Value = 45187
$ws1.Cells.Item(3798, 2).Value = 26741.46111094895
$ws1.Cells.Item(3799, 1).Value = 45188
$ws1.Cells.Item(3799, 2).Value = 27219.29687463422
$ws1.Cells.Item(3800, 1).Value = 45189
$ws1.Cells.Item(3800, 2).Value = 27115.84644697082
$ws1.Cells.Item(3801, 1).Value = 45190
$ws1.Cells.Item(3801, 2).Value = 26561.13345419872
$ws1.Cells.Item(3802, 1).Value = 45191
$ws1.Cells.Item(3802, 2).Value = 26572.0381115526
$ws1.Cells.Item(3803, 1).Value = 45192
$ws1.Cells.Item(3803, 2).Value = 26573.9234797301
$ws1.Cells.Item(3804, 1).Value = 45193
$ws1.Cells.Item(3804, 2).Value = 26249.5628978452
$ws1.Cells.Item(3805, 1).Value = 45194
$ws1.Cells.Item(3805, 2).Value = 26298.63467828123
$ws1.Cells.Item(3806, 1).Value = 45195
$ws1.Cells.Item(3806, 2).Value = 26204.75759083597
$ws1.Cells.Item(3807, 1).Value = 45196
$ws1.Cells.Item(3807, 2).Value = 26350.14689542806
$ws1.Cells.Item(3808, 1).Value = 45197
$ws1.Cells.Item(3808, 2).Value = 27009.01375072488
$ws1.Cells.Item(3809, 1).Value = 45198
$ws1.Cells.Item(3809, 2).Value = 26917.19910163798
$ws1.Cells.Item(3810, 1).Value = 45199
$ws1.Cells.Item(3810, 2).Value = 26969.87614407258
$ws1.Cells.Item(3811, 1).Value = 45200
$ws1.Cells.Item(3811, 2).Value = 27967.51057908711
$ws1.Cells.Item(3812, 1).Value = 45201
$ws1.Cells.Item(3812, 2).Value = 27615.06488509527
$ws1.Cells.Item(3813, 1).Value = 45202
$ws1.Cells.Item(3813, 2).Value = 27439.12194670512
$ws1.Cells.Item(3814, 1).Value = 45203
$ws1.Cells.Item(3814, 2).Value = 27792.1112772493
$ws1.Cells.Item(3815, 1).Value = 45204
$ws1.Cells.Item(3815, 2).Value = 27435.8746151351
$ws1.Cells.Item(3816, 1).Value = 45205
$ws1.Cells.Item(3816, 2).Value = 27958.19643735049
$ws1.Cells.Item(3817, 1).Value = 45206
$ws1.Cells.Item(3817, 2).Value = 27977.54349070001
$ws1.Cells.Item(3818, 1).Value = 45207
$ws1.Cells.Item(3818, 2).Value = 27948.10365174851
$ws1.Cells.Item(3819, 1).Value = 45208
$ws1.Cells.Item(3819, 2).Value = 27593.78253443967
$ws1.Cells.Item(3820, 1).Value = 45209
$ws1.Cells.Item(3820, 2).Value = 27392.2477027325
$ws1.Cells.Item(3821, 1).Value = 45210
$ws1.Cells.Item(3821, 2).Value = 26842.1904390924
$ws1.Cells.Item(3822, 1).Value = 45211
$ws1.Cells.Item(3822, 2).Value = 26729.13720581511
$ws1.Cells.Item(3823, 1).Value = 45212
$ws1.Cells.Item(3823, 2).Value = 26841.13622064441
$ws1.Cells.Item(3824, 1).Value = 45213
$ws1.Cells.Item(3824, 2).Value = 26863.18356907992
$ws1.Cells.Item(3825, 1).Value = 45214
$ws1.Cells.Item(3825, 2).Value = 27150.29700140705
$ws1.Cells.Item(3826, 1).Value = 45215
$ws1.Cells.Item(3826, 2).Value = 28513.30993247735
$ws1.Cells.Item(3827, 1).Value = 45216
$ws1.Cells.Item(3827, 2).Value = 28417.72175169982
$ws1.Cells.Item(3828, 1).Value = 45217
$ws1.Cells.Item(3828, 2).Value = 28328.24519813482
$ws1.Cells.Item(3829, 1).Value = 45218
$ws1.Cells.Item(3829, 2).Value = 28715.74814240795
$ws1.Cells.Item(3830, 1).Value = 45219
$ws1.Cells.Item(3830, 2).Value = 29677.39288847621
$ws1.Cells.Item(3831, 1).Value = 45220
$ws1.Cells.Item(3831, 2).Value = 29920.07449265145
$ws1.Cells.Item(3832, 1).Value = 45221
$ws1.Cells.Item(3832, 2).Value = 30019.38050086351
$ws1.Cells.Item(3833, 1).Value = 45222
$ws1.Cells.Item(3833, 2).Value = 32953.26276083098
$ws1.Cells.Item(3834, 1).Value = 45223
$ws1.Cells.Item(3834, 2).Value = 33846.72425733224
$ws1.Cells.Item(3835, 1).Value = 45224
$ws1.Cells.Item(3835, 2).Value = 34471.98603167202
$ws1.Cells.Item(3836, 1).Value = 45225
$ws1.Cells.Item(3836, 2).Value = 34174.45155291259
$ws1.Cells.Item(3837, 1).Value = 45226
$ws1.Cells.Item(3837, 2).Value = 33899.09305644032
$ws1.Cells.Item(3838, 1).Value = 45227
$ws1.Cells.Item(3838, 2).Value = 34092.63093283858
$ws1.Cells.Item(3839, 1).Value = 45228
$ws1.Cells.Item(3839, 2).Value = 34556.24281476162
$ws1.Cells.Item(3840, 1).Value = 45229
$ws1.Cells.Item(3840, 2).Value = 34498.70391946407
$ws1.Cells.Item(3841, 1).Value = 45230
$ws1.Cells.Item(3841, 2).Value = 34672.2892841885
$ws1.Cells.Item(3842, 1).Value = 45231
$ws1.Cells.Item(3842, 2).Value = 35457.45491210553
$ws1.Cells.Item(3843, 1).Value = 45232
$ws1.Cells.Item(3843, 2).Value = 34924.05545044328
$ws1.Cells.Item(3844, 1).Value = 45233
$ws1.Cells.Item(3844, 2).Value = 34731.38136896784
$ws1.Cells.Item(3845, 1).Value = 45234
$ws1.Cells.Item(3845, 2).Value = 35048.40783490107
$ws1.Cells.Item(3846, 1).Value = 45235
$ws1.Cells.Item(3846, 2).Value = 35061.92874919579
$ws1.Cells.Item(3847, 1).Value = 45236
$ws1.Cells.Item(3847, 2).Value = 35031.26888208706
$ws1.Cells.Item(3848, 1).Value = 45237
$ws1.Cells.Item(3848, 2).Value = 35436.53762957962
$ws1.Cells.Item(3849, 1).Value = 45238
$ws1.Cells.Item(3849, 2).Value = 35795.0806307102
$ws1.Cells.Item(3850, 1).Value = 45239
$ws1.Cells.Item(3850, 2).Value = 36768.42081912672
$ws1.Cells.Item(3851, 1).Value = 45240
$ws1.Cells.Item(3851, 2).Value = 37344.24900072035
$ws1.Cells.Item(3852, 1).Value = 45241
$ws1.Cells.Item(3852, 2).Value = 37122.72282430655
$ws1.Cells.Item(3853, 1).Value = 45242
$ws1.Cells.Item(3853, 2).Value = 37067.69698212008
$ws1.Cells.Item(3854, 1).Value = 45243
$ws1.Cells.Item(3854, 2).Value = 36549.16204829837
$ws1.Cells.Item(3855, 1).Value = 45244
$ws1.Cells.Item(3855, 2).Value = 35545.20143345407
$ws1.Cells.Item(3856, 1).Value = 45245
$ws1.Cells.Item(3856, 2).Value = 37903.66245166294
$ws1.Cells.Item(3857, 1).Value = 45246
$ws1.Cells.Item(3857, 2).Value = 36201.51611146142
$ws1.Cells.Item(3858, 1).Value = 45247
$ws1.Cells.Item(3858, 2).Value = 36527.76022530742
$ws1.Cells.Item(3859, 1).Value = 45248
$ws1.Cells.Item(3859, 2).Value = 36582.36844192274
$ws1.Cells.Item(3860, 1).Value = 45249
$ws1.Cells.Item(3860, 2).Value = 37413.99460790531
$ws1.Cells.Item(3861, 1).Value = 45250
$ws1.Cells.Item(3861, 2).Value = 37489.29847080202
$ws1.Cells.Item(3862, 1).Value = 45251
$ws1.Cells.Item(3862, 2).Value = 35965.36036068078
$ws1.Cells.Item(3863, 1).Value = 45252
$ws1.Cells.Item(3863, 2).Value = 37464.83293220907
$ws1.Cells.Item(3864, 1).Value = 45253
$ws1.Cells.Item(3864, 2).Value = 37293.31612742673
$ws1.Cells.Item(3865, 1).Value = 45254
$ws1.Cells.Item(3865, 2).Value = 37738.93169747125
$ws1.Cells.Item(3866, 1).Value = 45255
$ws1.Cells.Item(3866, 2).Value = 37809.85286625321
$ws1.Cells.Item(3867, 1).Value = 45256
$ws1.Cells.Item(3867, 2).Value = 37491.83818600814
$ws1.Cells.Item(3868, 1).Value = 45257
$ws1.Cells.Item(3868, 2).Value = 37250.16905148115
$ws1.Cells.Item(3869, 1).Value = 45258
$ws1.Cells.Item(3869, 2).Value = 37802.23604377473

# --- Sheet2: SeriesInfo ---
$ws2 = $wb.Worksheets.Item("SeriesInfo")
$ws2.Range("B1").Value = 0
$ws2.Range("B4").Value = " Bitcoin (USD)"
$ws2.Range("B5").Value = "Trace_2"
$ws2.Range("A6").Value = "Source"
$ws2.Range("B6").Value = "coingecko"

# --- Remove OtherInfo sheet (sheet3) entirely ---
$ws3 = $wb.Worksheets.Item("OtherInfo")
$null = $ws3.Delete()

# Restore the originally active sheet/tab
$ws1.Activate()

Write-Host "Edit complete"
